$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Log")

# --- Row 38: first result of the "MIX" / "90" test group just logged ---
$ws.Range("C38").Value = 0.75390625
$ws.Range("C37").Copy()
$ws.Range("C38").PasteSpecial(-4122)   # xlPasteFormats - match C35:C37's percent style

$ws.Range("G38").Value = 42542
$ws.Range("G37").Copy()
$ws.Range("G38").PasteSpecial(-4122)   # xlPasteFormats - match G35:G37's date style

# --- Row 39 ---
$ws.Range("C39").Value = 0.69921875
$ws.Range("C37").Copy()
$ws.Range("C39").PasteSpecial(-4122)

$ws.Range("G39").Value = 42542
$ws.Range("G37").Copy()
$ws.Range("G39").PasteSpecial(-4122)

# --- Row 40 ---
$ws.Range("C40").Value = 0.74609375
$ws.Range("C37").Copy()
$ws.Range("C40").PasteSpecial(-4122)

$ws.Range("G40").Value = 42542
$ws.Range("G37").Copy()
$ws.Range("G40").PasteSpecial(-4122)

$excel.CutCopyMode = $false

# Recalculate the dependent E (Change) / F (Group Change) formulas.
$null = $excel.Calculate()

# Reflect the scrolled/selected state left by finishing data entry at row 40/41.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 17
$excel.ActiveWindow.ScrollColumn = 1
$null = $ws.Range("C41").Select()
